$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the tuning parameters based on test-riding results
$ws.Range("M5").Value = 0.4
$ws.Range("N8").Value = 18000

# Update the active selection to reflect where the author was working
$ws.Range("M5").Select()
